# Update "想去人数" (want-to-go count) figures in column F for rows 2-6
# on both the "展览" and "全部类型" worksheets, matching the latest
# gh-pages generated data output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 5934
    $ws.Range("F3").Value = 19
    $ws.Range("F4").Value = 177
    $ws.Range("F5").Value = 987
    $ws.Range("F6").Value = 82
}
